$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.615.46'
$ws.Range('E2').Value = '  -4.05%  '

$ws.Range('D3').Value = '2.977.79'
$ws.Range('E3').Value = '  -4.95%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '539.88'
$ws.Range('E5').Value = '  -5.52%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.65'
$ws.Range('E6').Value = '  -7.83%  '

$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('E8').Value = '  -0.83%  '

$ws.Range('D9').Value = '2.985.72'
$ws.Range('E9').Value = '  -5.04%  '

$ws.Range('E10').Value = '  -3.38%  '

$ws.Range('E11').Value = '  -6.93%  '

$ws.Range('E12').Value = '  -4.56%  '

$ws.Range('D13').Value = '3.495.55'
$ws.Range('E13').Value = '  -5.17%  '

$ws.Range('E14').Value = '  -1.70%  '

$ws.Range('D15').Value = '61.677.97'
$ws.Range('E15').Value = '  -3.99%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.65'
$ws.Range('E16').Value = '  -5.42%  '

$ws.Range('D17').Value = '2.982.86'
$ws.Range('E17').Value = '  -5.05%  '

$ws.Range('E18').Value = '  -5.24%  '

$ws.Range('E19').Value = '  -1.77%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.05'
$ws.Range('E20').Value = '  -3.51%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.50'
$ws.Range('E21').Value = '  -5.36%  '

$ws.Range('E22').Value = '  -5.24%  '

$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('E24').Value = '  -3.43%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.88'
$ws.Range('E25').Value = '  -4.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.472'
$ws.Range('E26').Value = '  -2.62%  '

$ws.Range('D27').Value = '3.099.09'
$ws.Range('E27').Value = '  -5.40%  '

$ws.Range('E28').Value = '  -3.51%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.09%  '

$ws.Range('D30').Value = '0.0₃0938'
$ws.Range('E30').Value = '  -7.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.18'
$ws.Range('E31').Value = '  -6.99%  '

$ws.Range('E33').Value = '  -4.54%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '20.48'
$ws.Range('E34').Value = '  -3.35%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '159.01'
$ws.Range('E35').Value = '  -1.53%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.61'
$ws.Range('E36').Value = '  -5.04%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.89'
$ws.Range('E37').Value = '  -5.87%  '

$ws.Range('E38').Value = '  -3.58%  '

$ws.Range('E39').Value = '  -5.60%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.54'
$ws.Range('E40').Value = '  -7.62%  '

$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.54'
$ws.Range('E41').Value = '  -2.43%  '

$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.91'
$ws.Range('E42').Value = '  -3.88%  '

$ws.Range('D43').Value = '2.412.20'
$ws.Range('E43').Value = '  -8.84%  '

$ws.Range('E44').Value = '  -6.42%  '

$ws.Range('E45').Value = '  -2.36%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0591'
$ws.Range('E46').Value = '  -3.29%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.08'
$ws.Range('E47').Value = '  -5.32%  '

$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.997'
$ws.Range('E48').Value = '  +0.01%  '

$ws.Range('E49').Value = '  -3.21%  '

$ws.Range('E50').Value = '  -2.40%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.79'
$ws.Range('E51').Value = '  -6.12%  '
